$d = $word.ActiveDocument

# --- Edit 1: Deliverables paragraph -- append a textWrapping break + new sentence ---
# Locate the paragraph by its distinctive leading text.
$deliverablesPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("The main component of the deliverables")) {
        $deliverablesPara = $p
        break
    }
}
if ($deliverablesPara -eq $null) {
    throw "Could not find the Deliverables paragraph"
}

$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The main component of the deliverables will be a poster outlining the new GRAPH/Z’s strengths and weaknesses and highlighting the changes that we made from the original project. The poster will also include the abstract and writeup needed for entering it into the Supercomputing conference. Deliverables will also include the finished GRAPH/Z processing system and information comparing it to GraphLab and other existing similar tools. Also included with be data from profiling and traces.</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">We will test the newly proposed partitioning algorithm on a single node on an instance of AWS EC2 using the dataset from Stanford Network Analysis Project.</w:t></w:r></w:p>'
[void]$deliverablesPara.Range.InsertXML($xml1)

# --- Edit 2: Conclusion paragraph -- append a new closing sentence ---
$conclusionPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("The original GRAPH/Z was underperforming")) {
        $conclusionPara = $p
        break
    }
}
if ($conclusionPara -eq $null) {
    throw "Could not find the Conclusion body paragraph"
}

$concRange = $conclusionPara.Range
$concRange.Collapse(0)
$concRange.InsertAfter(" By using the newly proposed partitioning algorithm on GRAPH/Z, we hope to achieve better performance of the GRAPH/Z system.")
